$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'308.83"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'-2.82%"
$ws.Range("E2").Style = "Normal"
$ws.Range("D3").Value = "'37.30"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "'-6.59%"
$ws.Range("E3").Style = "Normal"
$ws.Range("D4").Value = "'5.108"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "'-0.55%"
$ws.Range("E4").Style = "Normal"
$ws.Range("D5").Value = "'0.07858"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "'-4.31%"
$ws.Range("E5").Style = "Normal"
$ws.Range("D6").Value = "'1.959"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "'-5.33%"
$ws.Range("E6").Style = "Normal"
$ws.Range("D7").Value = "'4.381"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "'2.03%"
$ws.Range("E7").Style = "Normal"
$ws.Range("D8").Value = "'8.279"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "'-0.42%"
$ws.Range("E8").Style = "Normal"
$ws.Range("E9").Value = "'-6.56%"
$ws.Range("E9").Style = "Normal"
$ws.Range("D10").Value = "'0.9236"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "'-1.21%"
$ws.Range("E10").Style = "Normal"
$ws.Range("D11").Value = "'0.1334"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "'-2.12%"
$ws.Range("E11").Style = "Normal"
$ws.Range("D12").Value = "'0.1948"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "'-2.20%"
$ws.Range("E12").Style = "Normal"
$ws.Range("D13").Value = "'0.08937"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "'-1.10%"
$ws.Range("E13").Style = "Normal"
$ws.Range("D14").Value = "'0.03445"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "'-0.95%"
$ws.Range("E14").Style = "Normal"
$ws.Range("D15").Value = "'0.09702"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "'-1.07%"
$ws.Range("E15").Style = "Normal"
$ws.Range("D16").Value = "'0.001378"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "'-1.46%"
$ws.Range("E16").Style = "Normal"
$ws.Range("D17").Value = "'0.005955"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "'-5.21%"
$ws.Range("E17").Style = "Normal"
$ws.Range("E19").Value = "'-2.01%"
$ws.Range("E19").Style = "Normal"
$ws.Range("E20").Value = "'0.23%"
$ws.Range("E20").Style = "Normal"
$ws.Range("E22").Value = "'1.65%"
$ws.Range("E22").Style = "Normal"
$ws.Range("D23").Value = "'0.02107"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "'5,172.09%"
$ws.Range("E23").Style = "Normal"
$ws.Range("D24").Value = "'0.04334"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "'0.22%"
$ws.Range("E24").Style = "Normal"
$ws.Range("E25").Value = "'-0.67%"
$ws.Range("E25").Style = "Normal"
$ws.Range("D26").Value = "'0.004527"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "'-4.75%"
$ws.Range("E26").Style = "Normal"
$ws.Range("E27").Value = "'3.91%"
$ws.Range("E27").Style = "Normal"
$ws.Range("D39").Value = "'0.02278"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "'2.27%"
$ws.Range("E39").Style = "Normal"
$ws.Range("D40").Value = "'0.05016"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "'-3.94%"
$ws.Range("E40").Style = "Normal"
$ws.Range("D41").Value = "'0.007646"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "'-0.36%"
$ws.Range("E41").Style = "Normal"
$ws.Range("D42").Value = "'0.009825"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "'1.68%"
$ws.Range("E42").Style = "Normal"
$ws.Range("E43").Value = "'-2.39%"
$ws.Range("E43").Style = "Normal"
$ws.Range("D44").Value = "'0.002063"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "'-1.55%"
$ws.Range("E44").Style = "Normal"
$ws.Range("D45").Value = "'0.008442"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "'-8.22%"
$ws.Range("E45").Style = "Normal"
$ws.Range("D46").Value = "'0.00006764"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "'3.30%"
$ws.Range("E46").Style = "Normal"
$ws.Range("D47").Value = "'0.00000000751"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "'0.04%"
$ws.Range("E47").Style = "Normal"
$ws.Range("D48").Value = "'0.003005"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "'8.13%"
$ws.Range("E48").Style = "Normal"
$ws.Range("D49").Value = "'0.001302"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "'8.48%"
$ws.Range("E49").Style = "Normal"
$ws.Range("D50").Value = "'0.00002103"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "'0.04%"
$ws.Range("E50").Style = "Normal"
$ws.Range("D51").Value = "'0.0002002"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "'0.04%"
$ws.Range("E51").Style = "Normal"
